# Auto-generated edit script: applies the Tiamat_Profits.xlsx value changes
# across sheets ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 31114.285
$ws.Range("I62").Value = 41560
$ws.Range("K62").Value = 41560
$ws.Range("M62").Value = -40936
$ws.Range("H65").Value = 31114.285
$ws.Range("I65").Value = 41560
$ws.Range("K65").Value = 207800
$ws.Range("M65").Value = -204680
$ws.Range("H93").Value = 28000
$ws.Range("J93").Value = 28000
$ws.Range("L93").Value = 28000
$ws.Range("N93").Value = -32992
$ws.Range("H125").Value = 55556310
$ws.Range("I125").Value = 58824270
$ws.Range("J125").Value = 1000
$ws.Range("K125").Value = 529418430
$ws.Range("L125").Value = 9000
$ws.Range("M125").Value = -529415970
$ws.Range("N125").Value = -13920

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 11666.667
$ws.Range("I28").Value = 11666.667
$ws.Range("K28").Value = 11666.667
$ws.Range("M28").Value = -11474.667
$ws.Range("H70").Value = 90577
$ws.Range("J70").Value = 90577
$ws.Range("L70").Value = 90577
$ws.Range("N70").Value = -91117
$ws.Range("H73").Value = 90577
$ws.Range("J73").Value = 90577
$ws.Range("L73").Value = 90577
$ws.Range("N73").Value = -92449
$ws.Range("H99").Value = 11666.667
$ws.Range("I99").Value = 11666.667
$ws.Range("K99").Value = 11666.667
$ws.Range("M99").Value = -8671.666999999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H92").Value = 53600.25
$ws.Range("J92").Value = 53600.25
$ws.Range("L92").Value = 53600.25
$ws.Range("N92").Value = -58592.25

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29630.127
$ws.Range("I31").Value = 38392.938
$ws.Range("J31").Value = 10936.134
$ws.Range("K31").Value = 38392.938
$ws.Range("L31").Value = 10936.134
$ws.Range("M31").Value = -38097.938
$ws.Range("N31").Value = -11526.134
$ws.Range("H34").Value = 29630.127
$ws.Range("I34").Value = 38392.938
$ws.Range("J34").Value = 10936.134
$ws.Range("K34").Value = 38392.938
$ws.Range("L34").Value = 10936.134
$ws.Range("M34").Value = -38190.938
$ws.Range("N34").Value = -11340.134

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H33").Value = 2511.4
$ws.Range("I33").Value = 619.0833
$ws.Range("J33").Value = 10080.667
$ws.Range("K33").Value = 3714.4998
$ws.Range("L33").Value = 60484.00199999999
$ws.Range("M33").Value = -3431.4998
$ws.Range("N33").Value = -61050.00199999999
$ws.Range("H34").Value = 31250126
$ws.Range("J34").Value = 35714416
$ws.Range("L34").Value = 107143248
$ws.Range("N34").Value = -107143416
$ws.Range("H39").Value = 111113770
$ws.Range("J39").Value = 111113770
$ws.Range("L39").Value = 333341310
$ws.Range("N39").Value = -333341898
$ws.Range("H55").Value = 27087440
$ws.Range("I55").Value = 1700
$ws.Range("J55").Value = 33337994
$ws.Range("K55").Value = 5100
$ws.Range("L55").Value = 100013982
$ws.Range("M55").Value = -4923
$ws.Range("N55").Value = -100014336
$ws.Range("H64").Value = 2334835.5
$ws.Range("I64").Value = 1500
$ws.Range("J64").Value = 2759078.2
$ws.Range("K64").Value = 4500
$ws.Range("L64").Value = 8277234.600000001
$ws.Range("M64").Value = -4230
$ws.Range("N64").Value = -8277774.600000001
$ws.Range("H67").Value = 2334835.5
$ws.Range("I67").Value = 1500
$ws.Range("J67").Value = 2759078.2
$ws.Range("K67").Value = 4500
$ws.Range("L67").Value = 8277234.600000001
$ws.Range("M67").Value = -3564
$ws.Range("N67").Value = -8279106.600000001
$ws.Range("H80").Value = 2355.889
$ws.Range("J80").Value = 2883.8333
$ws.Range("L80").Value = 8651.499899999999
$ws.Range("N80").Value = -10523.4999
$ws.Range("H83").Value = 2355.889
$ws.Range("J83").Value = 2883.8333
$ws.Range("L83").Value = 25954.4997
$ws.Range("N83").Value = -35314.4997
$ws.Range("H125").Value = 8507.0625
$ws.Range("J125").Value = 8607.532999999999
$ws.Range("L125").Value = 25822.599
$ws.Range("N125").Value = -35662.599
$ws.Range("H130").Value = 3499
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 3499
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 10497
$ws.Range("M130").ClearContents()
$ws.Range("N130").Value = -20537

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1322.95
$ws.Range("I97").Value = 1241.5
$ws.Range("K97").Value = 1241.5
$ws.Range("M97").Value = -745.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1033.3334
$ws.Range("I46").Value = 600
$ws.Range("J46").Value = 1250
$ws.Range("K46").Value = 600
$ws.Range("L46").Value = 1250
$ws.Range("M46").Value = -412
$ws.Range("N46").Value = -1626
$ws.Range("H76").Value = 24333
$ws.Range("J76").Value = 24333
$ws.Range("L76").Value = 24333
$ws.Range("N76").Value = -25009
$ws.Range("H79").Value = 24333
$ws.Range("J79").Value = 24333
$ws.Range("L79").Value = 24333
$ws.Range("N79").Value = -26673

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H63").Value = 25740.666
$ws.Range("J63").Value = 25740.666
$ws.Range("L63").Value = 25740.666
$ws.Range("N63").Value = -26988.666
$ws.Range("H66").Value = 25740.666
$ws.Range("J66").Value = 25740.666
$ws.Range("L66").Value = 77221.99800000001
$ws.Range("N66").Value = -83461.99800000001
$ws.Range("H69").Value = 10500
$ws.Range("J69").Value = 10500
$ws.Range("L69").Value = 10500
$ws.Range("N69").Value = -11998
$ws.Range("H72").Value = 10500
$ws.Range("J72").Value = 10500
$ws.Range("L72").Value = 31500
$ws.Range("N72").Value = -38988
$ws.Range("H80").Value = 42450.25
$ws.Range("J80").Value = 42450.25
$ws.Range("L80").Value = 42450.25
$ws.Range("N80").Value = -44446.25
$ws.Range("H83").Value = 42450.25
$ws.Range("J83").Value = 42450.25
$ws.Range("L83").Value = 127350.75
$ws.Range("N83").Value = -137334.75
$ws.Range("H92").Value = 24000
$ws.Range("J92").Value = 24000
$ws.Range("L92").Value = 24000
$ws.Range("N92").Value = -28992
$ws.Range("H126").Value = 844.2
$ws.Range("I126").Value = 602.1177
$ws.Range("K126").Value = 1806.3531
$ws.Range("M126").Value = 663.6469
